$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): update view count F14 for 苏州·COME IN JOY event
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F14").Value = 12733

# Sheet "演出" (Shows): update view count F3 for 苏州·足太Penta生日会2024 event
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 96

# Sheet "全部类型" (All types): mirrors the same two events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F15").Value = 12733
$wsAll.Range("F16").Value = 96
